$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 234, shifting existing rows 234-280 down to 235-281
$ws.Rows.Item(234).Insert()

# Populate the newly inserted row 234 with the new weekly record
$ws.Cells.Item(234, 1).Value = 11
$ws.Cells.Item(234, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(234, 3).Value = "Bíobío"
$ws.Cells.Item(234, 4).Value = 45211
$ws.Cells.Item(234, 5).Value = 8
$ws.Cells.Item(234, 6).Value = 100112032
$ws.Cells.Item(234, 7).Value = "Zapallo italiano"
$ws.Cells.Item(234, 8).Value = "Sin especificar"
$ws.Cells.Item(234, 9).Value = "Primera"
$ws.Cells.Item(234, 10).Value = 120
$ws.Cells.Item(234, 11).Value = 20000
$ws.Cells.Item(234, 12).Value = 20000
$ws.Cells.Item(234, 13).Value = 20000
$ws.Cells.Item(234, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(234, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(234, 16).Value = 400
$ws.Cells.Item(234, 17).Value = 50
$ws.Cells.Item(234, 18).Value = "Hortaliza"
